# ECGToolKit ReadMe.docx edit script
# Applies the changes described by the target diff:
#  1. Inserts a new preamble block (license blurb, attribution, repo link,
#     blank line, bold "USER MANUAL" heading, blank line) before the
#     existing first paragraph.
#  2. Cleans up several runs that used to be split around proofing marks
#     (spell-check / grammar-check squiggles) so the visible text reads the
#     same but is simplified.
#  3. Appends " from" to the "Launch the ECGViewer" bullet.
#  4. Splits "Load the Plugin for Muse" into "...for M" / "USE" and moves the
#     "_GoBack" bookmark from its old location (after the first "exported
#     file" sentence) to just after that new "USE" run.
#  5. Removes the stray paragraph-mark rFonts hint on an empty bullet.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. New preamble block, inserted before the very first paragraph.
# ---------------------------------------------------------------------
$first = $d.Paragraphs(1).Range
$first.InsertBefore("This is open-source software Licensed under the Apache License, Version 2.0`r" + `
    "The new functions were developed on C#ECKtoolkit 2.4 initiated by Maarten van Ettinger`r" + `
    "https://git.code.sf.net/p/ecgtoolkit-cs/git ecgtoolkit-cs-git`r" + `
    "`r" + `
    "USER MANUAL`r" + `
    "`r")

# Bold the "USER MANUAL" heading (paragraph 5 of the new block).
$manual = $d.Paragraphs(5).Range
$manual.Font.Bold = 1
$manual.Font.BoldBi = 1

Write-Host "p1:" $d.Paragraphs(1).Range.Text
Write-Host "p5:" $d.Paragraphs(5).Range.Text
Write-Host "p7:" $d.Paragraphs(7).Range.Text
